$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the shared formulas in columns G (Purple CFU/mL) and I (Relative
# Abundance) down through the newly-populated rows 57:61 so every row in
# the 33:61 block is consistent (matches H and M, which already run through
# row 61).
for ($r = 57; $r -le 61; $r++) {
    $ws.Range("G$r").Formula = "=(20*(1/D$r))*E$r"
    $ws.Range("G$r").Font.Name = "Times"
    $ws.Range("I$r").Formula = "=G$r/(G$r+H$r)"
    $ws.Range("I$r").Font.Name = "Times"
}

# Move the active selection off the old E56:E61 block onto H31, matching
# where the author was last working.
$ws.Range("H31").Select() | Out-Null
